# "add Crisis and Credit Allocation" -- update the regression coefficient
# table on Sheet1 with the new (2005-12-31 to 2008-09-30) estimates.
#
# Each target cell currently stores a number-looking value as literal TEXT
# (shared string), not a numeric cell. Assigning a plain numeric-looking
# string straight to .Value would get auto-coerced into a real number, so
# for every cell we briefly force Text number-formatting, assign the new
# text, then restore the cell's original ("Normal") style so the stored
# formatting stays exactly as it was before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "B2" "0.17"
Set-TextValue "B3" "-0.01"
Set-TextValue "B4" "-0.09"

Set-TextValue "C2" "44.29***"
Set-TextValue "C3" "2.21***"
Set-TextValue "C4" "0.98"

Set-TextValue "D2" "-0.89"
Set-TextValue "D3" "0.46***"
Set-TextValue "D4" "0.82*"
